$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: snapshot formatting of current last row (403) for reuse on new last row (413) ---
$ws.Range("A403:N403").Copy()
$ws.Range("A413:N413").PasteSpecial(-4122)

# --- Step 2: re-style row 403 as a normal (non-last) row, matching template row 400 ---
$ws.Range("A400:N400").Copy()
$ws.Range("A403:N403").PasteSpecial(-4122)

# --- Step 3: write new rows 404-412 with values + banded-row formatting ---
$ws.Range("A400:N400").Copy()
$ws.Range("A404:N404").PasteSpecial(-4122)
$ws.Range("A404").Value = 45610.96508303241
$ws.Range("B404").Value = 'tngusvhs@gmail.com'
$ws.Range("C404").Value = '생명과학과'
$ws.Range("D404").Value = 20243529
$ws.Range("E404").Value = '이수현'
$ws.Range("F404").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G404").Value = '대공황, 대번영, 대침체'
$ws.Range("H404").Value = '"60%"'
$ws.Range("I404").Value = '"45%"'
$ws.Range("J404").Value = '황금기에 소득증가율이 가장 높은 계층은 하위 20%(1분위)이다.'
$ws.Range("K404").Value = '천지비'
$ws.Range("L404").Value = 'Black'
$ws.Range("N404").Value = 'A : 100% 확률로 900불 잃기'
$ws.Range("M404").Clear()

$ws.Range("A401:N401").Copy()
$ws.Range("A405:N405").PasteSpecial(-4122)
$ws.Range("A405").Value = 45610.984333229164
$ws.Range("B405").Value = 'jerryterryharry@gmail.com'
$ws.Range("C405").Value = '빅데이터'
$ws.Range("D405").Value = 20205162
$ws.Range("E405").Value = '문진영'
$ws.Range("F405").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G405").Value = '대공황, 대번영, 대침체'
$ws.Range("H405").Value = '"60%"'
$ws.Range("I405").Value = '"45%"'
$ws.Range("J405").Value = '황금기에 소득이 가장 많이 늘어난 계층은 하위 20%(1분위)이다.'
$ws.Range("K405").Value = '천지비'
$ws.Range("L405").Value = 'Red'
$ws.Range("M405").Value = 'A : 100% 확률로 900불 따기'
$ws.Range("N405").Clear()

$ws.Range("A400:N400").Copy()
$ws.Range("A406:N406").PasteSpecial(-4122)
$ws.Range("A406").Value = 45611.03941334491
$ws.Range("B406").Value = 'lhw2565@gmail.com'
$ws.Range("C406").Value = '미디어스쿨'
$ws.Range("D406").Value = 20242565
$ws.Range("E406").Value = '이혜원'
$ws.Range("F406").Value = '소득불평등이 심하면 건강 및 사회문제지수가 나빠진다.'
$ws.Range("G406").Value = '대번영, 대공황, 대침체'
$ws.Range("H406").Value = '"50%"'
$ws.Range("I406").Value = '"35%"'
$ws.Range("J406").Value = '황금기에 소득증가율이 가장 높은 계층은 하위 20%(1분위)이다.'
$ws.Range("K406").Value = '산택손'
$ws.Range("L406").Value = 'Black'
$ws.Range("N406").Value = 'B :  90% 확률로 1,000불 잃기'
$ws.Range("M406").Clear()

$ws.Range("A399:N399").Copy()
$ws.Range("A407:N407").PasteSpecial(-4122)
$ws.Range("A407").Value = 45611.464987430554
$ws.Range("B407").Value = 'yhh323@naver.com'
$ws.Range("C407").Value = '체육학과'
$ws.Range("D407").Value = 20184132
$ws.Range("E407").Value = '유형호'
$ws.Range("F407").Value = '소득이 많은 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G407").Value = '대침체, 대공황, 대번영'
$ws.Range("H407").Value = '"50%"'
$ws.Range("I407").Value = '"45%"'
$ws.Range("J407").Value = '황금기에 소득증가율이 가장 높은 계층은 하위 20%(1분위)이다.'
$ws.Range("K407").Value = '천지비'
$ws.Range("L407").Value = 'Black'
$ws.Range("N407").Value = 'B :  90% 확률로 1,000불 잃기'
$ws.Range("M407").Clear()

$ws.Range("A400:N400").Copy()
$ws.Range("A408:N408").PasteSpecial(-4122)
$ws.Range("A408").Value = 45611.540331122684
$ws.Range("B408").Value = 'bigeyejimmy1@naver.com'
$ws.Range("C408").Value = '경영학과'
$ws.Range("D408").Value = 20182850
$ws.Range("E408").Value = '김현준'
$ws.Range("F408").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G408").Value = '대공황, 대번영, 대침체'
$ws.Range("H408").Value = '"60%"'
$ws.Range("I408").Value = '"45%"'
$ws.Range("J408").Value = '신자유주의시기에 소득증가율이 가장 높은 계층은 하위 20%(1분위)이다.'
$ws.Range("K408").Value = '풍뢰익'
$ws.Range("L408").Value = 'Black'
$ws.Range("N408").Value = 'A : 100% 확률로 900불 잃기'
$ws.Range("M408").Clear()

$ws.Range("A401:N401").Copy()
$ws.Range("A409:N409").PasteSpecial(-4122)
$ws.Range("A409").Value = 45611.588682002315
$ws.Range("B409").Value = 'jiyewon5555@gmail.com'
$ws.Range("C409").Value = '광고홍보학과'
$ws.Range("D409").Value = 20202641
$ws.Range("E409").Value = '지예원'
$ws.Range("F409").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G409").Value = '대공황, 대번영, 대침체'
$ws.Range("H409").Value = '"60%"'
$ws.Range("I409").Value = '"35%"'
$ws.Range("J409").Value = '황금기에 소득증가율이 가장 높은 계층은 하위 20%(1분위)이다.'
$ws.Range("K409").Value = '천지비'
$ws.Range("L409").Value = 'Red'
$ws.Range("M409").Value = 'A : 100% 확률로 900불 따기'
$ws.Range("N409").Clear()

$ws.Range("A400:N400").Copy()
$ws.Range("A410:N410").PasteSpecial(-4122)
$ws.Range("A410").Value = 45611.59365976852
$ws.Range("B410").Value = 'a01053076907@gmail.com'
$ws.Range("C410").Value = '융합과학수사학과'
$ws.Range("D410").Value = 20246927
$ws.Range("E410").Value = '원은성'
$ws.Range("F410").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G410").Value = '대공황, 대침체, 대번영'
$ws.Range("H410").Value = '"50%"'
$ws.Range("I410").Value = '"45%"'
$ws.Range("J410").Value = '신자유주의시기에 소득이 가장 많이 늘어난 계층은 하위 20%(1분위)이다.'
$ws.Range("K410").Value = '풍뢰익'
$ws.Range("L410").Value = 'Black'
$ws.Range("N410").Value = 'B :  90% 확률로 1,000불 잃기'
$ws.Range("M410").Clear()

$ws.Range("A401:N401").Copy()
$ws.Range("A411:N411").PasteSpecial(-4122)
$ws.Range("A411").Value = 45611.65700553241
$ws.Range("B411").Value = 'hyj13223@naver.com'
$ws.Range("C411").Value = '정치행정학과'
$ws.Range("D411").Value = 20212432
$ws.Range("E411").Value = '이현진'
$ws.Range("F411").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G411").Value = '대공황, 대번영, 대침체'
$ws.Range("H411").Value = '"60%"'
$ws.Range("I411").Value = '"45%"'
$ws.Range("J411").Value = '황금기에 소득이 가장 많이 늘어난 계층은 하위 20%(1분위)이다.'
$ws.Range("K411").Value = '천지비'
$ws.Range("L411").Value = 'Red'
$ws.Range("M411").Value = 'A : 100% 확률로 900불 따기'
$ws.Range("N411").Clear()

$ws.Range("A400:N400").Copy()
$ws.Range("A412:N412").PasteSpecial(-4122)
$ws.Range("A412").Value = 45611.703416331016
$ws.Range("B412").Value = 'bcy1976@naver.com'
$ws.Range("C412").Value = '빅데이터학과'
$ws.Range("D412").Value = 20235180
$ws.Range("E412").Value = '변치윤'
$ws.Range("F412").Value = '소득불평등이 심한 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G412").Value = '대공황, 대번영, 대침체'
$ws.Range("H412").Value = '"50%"'
$ws.Range("I412").Value = '"45%"'
$ws.Range("J412").Value = '황금기에 소득증가율이 가장 높은 계층은 하위 20%(1분위)이다.'
$ws.Range("K412").Value = '천지비'
$ws.Range("L412").Value = 'Black'
$ws.Range("N412").Value = 'A : 100% 확률로 900불 잃기'
$ws.Range("M412").Clear()

# --- Step 4: write new last row 413 (formatting already copied from original row 403) ---
$ws.Range("A413").Value = 45611.709958182866
$ws.Range("B413").Value = 'emf1811@naver.com'
$ws.Range("C413").Value = '바이오메디컬학과'
$ws.Range("D413").Value = 20233605
$ws.Range("E413").Value = '김들'
$ws.Range("F413").Value = '소득이 많은 나라에서 건강 및 사회문제지수가 나쁘게 나온다.'
$ws.Range("G413").Value = '대공황, 대침체, 대번영'
$ws.Range("H413").Value = '"60%"'
$ws.Range("I413").Value = '"35%"'
$ws.Range("J413").Value = '황금기에 소득이 가장 많이 늘어난 계층은 하위 20%(1분위)이다.'
$ws.Range("K413").Value = '지천태'
$ws.Range("L413").Value = 'Black'
$ws.Range("N413").Value = 'B :  90% 확률로 1,000불 잃기'
$ws.Range("M413").Clear()

# --- Step 5: resize the table (ListObject) to include the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N413"))

# --- Step 6: update selection to match the author's final cursor position ---
$ws.Range("A420").Select()